$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row -> @(newPriceOrNull, newVolumeOrNull)
$updates = @{
    2  = @("308.98", "0.36%")
    3  = @("40.98", "-0.17%")
    4  = @("5.135", "1.83%")
    5  = @("0.07632", "-0.14%")
    6  = @("1.608", "-0.27%")
    8  = @($null, "0.11%")
    9  = @("0.1271", "26.03%")
    10 = @("0.1805", "2.07%")
    11 = @("0.09081", "-1.52%")
    12 = @("0.04338", "2.68%")
    13 = @($null, "-0.69%")
    14 = @("0.001250", "-0.47%")
    15 = @("0.005658", "-2.45%")
    16 = @("3.353", "-0.11%")
    17 = @("4.286", "0.60%")
    18 = @("0.3346", "2.35%")
    19 = @("6.897", "1.62%")
    20 = @("0.1383", "2.56%")
    21 = @("0.2735", "0.49%")
    22 = @($null, "-2.80%")
    23 = @("0.001270", "4.30%")
    24 = @("0.004062", "-0.54%")
    25 = @($null, "-2.16%")
    26 = @($null, "24.78%")
    38 = @("0.02423", "0.54%")
    39 = @("0.05226", "0.89%")
    40 = @("0.007858", "0.91%")
    41 = @("0.1301", "-0.52%")
    42 = @("0.006787", "-4.31%")
    43 = @($null, "-6.52%")
    44 = @("0.007446", "-0.40%")
    45 = @("0.3354", "9.80%")
    46 = @("0.00006869", "7.59%")
    47 = @("0.00000000750", "0.03%")
    48 = @("0.1625", "2,530.91%")
    49 = @($null, "-31.74%")
    50 = @("0.00002101", "0.03%")
    51 = @("0.0002000", "0.03%")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $priceVal = $pair[0]
    $volVal = $pair[1]

    if ($priceVal -ne $null) {
        $cell = $ws.Cells.Item($row, 4)
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $priceVal
        $cell.Style = $origStyle
    }

    if ($volVal -ne $null) {
        $cell = $ws.Cells.Item($row, 5)
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $volVal
        $cell.Style = $origStyle
    }
}
